# AliExpress 2.0 7th upload - reference sheet column rework
# Rearranges header row (row 1), drops unused labels, adds the new
# crawling columns (video link / shipping method / shipping fee / etc.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final header labels for row 1, written in the same order the columns
# were actually filled in (kept labels first where the text is unchanged,
# then the newly-typed labels: search/price columns, then the M-P block,
# then the J-K-L video/shipping block added last per the commit message).
$headers = [ordered]@{
    "A1" = "검색키워드"
    "B1" = "정가"
    "C1" = "할인가"
    "D1" = "댓글수"
    "E1" = "구매수"
    "F1" = "이미지저장경로"
    "G1" = "네이버카테고리"
    "H1" = "국내사이트제목"
    "I1" = "국내사이트태그"
    "M1" = "검색메인cate"
    "N1" = "상품url"
    "O1" = "상품가격범위"
    "P1" = "상품명"
    "Q1" = "대표이미지"
    "R1" = "옵션1"
    "S1" = "옵션2"
    "T1" = "옵션_종합"
    "U1" = "상세페이지"
    "V1" = "고객사_상품코드"
    "J1" = "동영상링크"
    "K1" = "배송방법"
    "L1" = "배송비"
}

foreach ($addr in $headers.Keys) {
    $ws.Range($addr).Value = $headers[$addr]
}

# The previous header row used the "표준 2" (vertical-center) cell style;
# the refreshed sheet goes back to the plain default style.
$ws.Range("A1:V1").Style = "표준"

# Selection moved to M1 when the sheet was last saved.
$ws.Range("M1").Select()
